$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column A, shifting existing columns (A->B ... K->L) to the right
$ws.Columns.Item(1).Insert()

# New header "Name" in A1, matching the style used by the other headers
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2 -- update to "PA1NT #Peak"
$ws.Range("A2").Value = "PA1NT #Peak"
$ws.Range("B2").Value = 1.18
$ws.Range("C2").Value = 47
$ws.Range("D2").Value = "Jett"
$ws.Range("E2").Value = 33.5
$ws.Range("F2").Value = 47
$ws.Range("G2").Value = 486
$ws.Range("H2").Value = 446
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 29
$ws.Range("K2").Value = "Lead"
$ws.Range("L2").Value = "['Engager', 'Sniper', 'Rusher']"

# Row 3 -- add new player "big bean #jakee"
$ws.Range("A3").Value = "big bean #jakee"
$ws.Range("B3").Value = 1.11
$ws.Range("C3").Value = 46
$ws.Range("D3").Value = "Clove"
$ws.Range("E3").Value = 32.7
$ws.Range("F3").Value = 64
$ws.Range("G3").Value = 579
$ws.Range("H3").Value = 557
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 101
$ws.Range("K3").Value = "Titanium"
$ws.Range("L3").Value = "['Clutch King', 'Engager', 'Sniper', 'Rusher']"

Write-Output "done"
